$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (some values look numeric,
# e.g. "211.96", and would otherwise be auto-converted to a Double by Excel).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.646.74"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "1.592.99"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "211.96"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.0617"
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("E9").Value = "  -2.85%  "
$ws.Range("D10").Value = "19.61"
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "1.817.13"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "1.596.27"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").Value = "65.24"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "26.612.39"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -2.62%  "
$ws.Range("D19").Value = "208.53"
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "6.70"
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").Value = "2.31"
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("D24").Value = "8.88"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "145.95"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "7.16"
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").Value = "15.31"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").Value = "0.666"
$ws.Range("E33").Value = "  -10.65%  "
$ws.Range("D34").Value = "2.90"
$ws.Range("D35").Value = "1.300.04"
$ws.Range("E35").Value = "  -3.38%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  -4.99%  "
$ws.Range("E38").Value = "  -4.04%  "
$ws.Range("D39").Value = "0.829"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").Value = "2.18"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "63.22"
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("D45").Value = "1.729.86"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").Value = "89.06"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "1.60"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").Value = "0.806"
$ws.Range("E48").Value = "  -7.35%  "
$ws.Range("D49").Value = "0.0982"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("E50").Value = "  -2.58%  "
$ws.Range("D51").Value = "7.52"
$ws.Range("E51").Value = "  -2.14%  "

# Reset column D cells back to the default (unstyled) cell style now that
# the text value is set, so no stray number-format style lingers on them.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
